# Daily attendance processing - reverse the order of the comma-separated
# entries in the "Recorded By" column (G) on every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value()

    if ($null -eq $val) { continue }
    if ($val -notmatch ",") { continue }

    $parts = $val -split ",\s*"
    $n = $parts.Count

    $reversed = @()
    for ($i = $n - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }

    $cell.Value = [string]::Join(", ", $reversed)
}
